$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header labels in row 3 (weighting headers for each review).
# "30% Week 5" -> "10% Week 5" and "10% Week 13" -> "30% Week 13"
$ws.Range("C3").Value = "10% Week 5"
$ws.Range("F3").Value = "30% Week 13"

# Fill in the remaining marks (Code Review 2, Code Review 3, Final Deliverable)
# for every contributor - each gets 25 marks, bringing each row's total to 100.
$ws.Range("D8:F11").Value = 25

# Update the selection to match the saved state (entire column F selected,
# active cell F1).
$ws.Range("F:F").Select()
